# October push - Update xMDA, add SOPH, etc
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variableNames")

# Insert a new row at row 22, shifting existing rows 22-26 down to 23-27
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the F37 / hispanicGroup variable info
$ws.Cells.Item(22, 3).Value = "F37"
$ws.Cells.Item(22, 5).Value = "F37"
$ws.Cells.Item(22, 1).Value = "hispanicGroup"
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 4).Value = 1

# Match the style (alignment/border/fill) used by column C/E in the other rows
$ws.Cells.Item(23, 3).Copy()
$ws.Cells.Item(22, 3).PasteSpecial(-4122)
$ws.Cells.Item(23, 5).Copy()
$ws.Cells.Item(22, 5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view to match the recorded selection
$ws.Range("A29").Select()
